$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 03:18:49"

$ws1.Range("A6").Value = "03:18:49"
$ws1.Range("B6").Value = "03:48"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 30
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = "03:18:49"
$ws1.Range("B7").Value = "04:01"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 43
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = "03:18:49"
$ws1.Range("B8").Value = "04:47"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 89
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = "03:18:49"
$ws1.Range("B9").Value = "04:53"
$ws1.Range("C9").Value = "11_ETCHEVERRY"
$ws1.Range("D9").Value = 95
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "03:18:49"
$ws1.Range("B10").Value = "05:16"
$ws1.Range("C10").Value = "17_ROMERO"
$ws1.Range("D10").Value = 118
$ws1.Range("E10").Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 03:18:49"
$ws2.Range("A3").Value = "Total filas: 0"
# Remove the header row (5) and the single data row (6); deleting row 5 twice
# shifts row 6 into row 5's place, removing both.
$ws2.Rows.Item(5).Delete()
$ws2.Rows.Item(5).Delete()

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 03:18:49"

Write-Host "edit complete"
